$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.325956666666666
$ws.Range("H2").Value = 18.97787
$ws.Range("I2").Value = 0.4468357575736242
$ws.Range("J2").Value = 0.4592138460625664
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 6.041227666666667
$ws.Range("N2").Value = 18.123683
$ws.Range("O2").Value = 0.1819046328309099
$ws.Range("P2").Value = 0.201743010396314
$ws.Range("Q2").Value = 38.21654443280111
$ws.Range("R2").Value = 343.94889989521
$ws.Range("S2").Value = 0.08128149441715157
$ws.Range("T2").Value = 0.09264318372033166

$ws.Range("G3").Value = 6.325956666666666
$ws.Range("H3").Value = 18.97787
$ws.Range("I3").Value = 0.4468357575736242
$ws.Range("J3").Value = 0.4592138460625664
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 9.354491666666666
$ws.Range("N3").Value = 28.063475
$ws.Range("O3").Value = 0.2816688040634135
$ws.Range("P3").Value = 0.3123873844340412
$ws.Range("Q3").Value = 59.17610892202777
$ws.Range("R3").Value = 532.58498029825
$ws.Range("S3").Value = 0.1258596934485321
$ws.Range("T3").Value = 0.1434526122673816

$ws.Range("G4").Value = 6.325956666666666
$ws.Range("H4").Value = 18.97787
$ws.Range("I4").Value = 0.4468357575736242
$ws.Range("J4").Value = 0.4592138460625664
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.096198666666667
$ws.Range("N4").Value = 15.288596
$ws.Range("O4").Value = 0.1534492984610312
$ws.Range("P4").Value = 0.1701843594247949
$ws.Range("Q4").Value = 32.23833193005778
$ws.Range("R4").Value = 290.14498737052
$ws.Range("S4").Value = 0.06856663352697605
$ws.Range("T4").Value = 0.07815101423115424

$ws.Range("G5").Value = 6.325956666666666
$ws.Range("H5").Value = 18.97787
$ws.Range("I5").Value = 0.4468357575736242
$ws.Range("J5").Value = 0.4592138460625664
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 9.797388999999999
$ws.Range("N5").Value = 19.594778
$ws.Range("O5").Value = 0.2950046823396649
$ws.Range("P5").Value = 0.2181184421382489
$ws.Range("Q5").Value = 61.97785826047666
$ws.Range("R5").Value = 371.8671495628599
$ws.Range("S5").Value = 0.1318186407210105
$ws.Range("T5").Value = 0.1001630087114806

$ws.Range("G6").Value = 6.325956666666666
$ws.Range("H6").Value = 18.97787
$ws.Range("I6").Value = 0.4468357575736242
$ws.Range("J6").Value = 0.4592138460625664
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 2.921654
$ws.Range("N6").Value = 8.764962000000001
$ws.Range("O6").Value = 0.08797258230498059
$ws.Range("P6").Value = 0.09756680360660122
$ws.Range("Q6").Value = 18.48225659899333
$ws.Range("R6").Value = 166.34030939094
$ws.Range("S6").Value = 0.03930929545995401
$ws.Range("T6").Value = 0.04480402713221842

$ws.Range("G7").Value = 6.686451000000001
$ws.Range("H7").Value = 20.059353
$ws.Range("I7").Value = 0.4722993778644153
$ws.Range("J7").Value = 0.4853828506917099
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 6.041227666666667
$ws.Range("N7").Value = 18.123683
$ws.Range("O7").Value = 0.1819046328309099
$ws.Range("P7").Value = 0.201743010396314
$ws.Range("Q7").Value = 40.39437277301101
$ws.Range("R7").Value = 363.549354957099
$ws.Range("S7").Value = 0.08591344491669364
$ws.Range("T7").Value = 0.09792259749329014

$ws.Range("G8").Value = 6.686451000000001
$ws.Range("H8").Value = 20.059353
$ws.Range("I8").Value = 0.4722993778644153
$ws.Range("J8").Value = 0.4853828506917099
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 9.354491666666666
$ws.Range("N8").Value = 28.063475
$ws.Range("O8").Value = 0.2816688040634135
$ws.Range("P8").Value = 0.3123873844340412
$ws.Range("Q8").Value = 62.54835015907501
$ws.Range("R8").Value = 562.9351514316751
$ws.Range("S8").Value = 0.1330320009229641
$ws.Range("T8").Value = 0.151627479176722

$ws.Range("G9").Value = 6.686451000000001
$ws.Range("H9").Value = 20.059353
$ws.Range("I9").Value = 0.4722993778644153
$ws.Range("J9").Value = 0.4853828506917099
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 5.096198666666667
$ws.Range("N9").Value = 15.288596
$ws.Range("O9").Value = 0.1534492984610312
$ws.Range("P9").Value = 0.1701843594247949
$ws.Range("Q9").Value = 34.075482670932
$ws.Range("R9").Value = 306.6793440383881
$ws.Range("S9").Value = 0.07247400819687602
$ws.Range("T9").Value = 0.08260456952074952

$ws.Range("G10").Value = 6.686451000000001
$ws.Range("H10").Value = 20.059353
$ws.Range("I10").Value = 0.4722993778644153
$ws.Range("J10").Value = 0.4853828506917099
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 9.797388999999999
$ws.Range("N10").Value = 19.594778
$ws.Range("O10").Value = 0.2950046823396649
$ws.Range("P10").Value = 0.2181184421382489
$ws.Range("Q10").Value = 65.509761476439
$ws.Range("R10").Value = 393.058568858634
$ws.Range("S10").Value = 0.1393305279361132
$ws.Range("T10").Value = 0.105870951233498

$ws.Range("G11").Value = 6.686451000000001
$ws.Range("H11").Value = 20.059353
$ws.Range("I11").Value = 0.4722993778644153
$ws.Range("J11").Value = 0.4853828506917099
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 2.921654
$ws.Range("N11").Value = 8.764962000000001
$ws.Range("O11").Value = 0.08797258230498059
$ws.Range("P11").Value = 0.09756680360660122
$ws.Range("Q11").Value = 19.535496309954
$ws.Range("R11").Value = 175.819466789586
$ws.Range("S11").Value = 0.04154939589176841
$ws.Range("T11").Value = 0.0473572532674503

$ws.Range("G12").Value = 1.1448225
$ws.Range("H12").Value = 2.289645
$ws.Range("I12").Value = 0.08086486456196039
$ws.Range("J12").Value = 0.05540330324572383
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 6.041227666666667
$ws.Range("N12").Value = 18.123683
$ws.Range("O12").Value = 0.1819046328309099
$ws.Range("P12").Value = 0.201743010396314
$ws.Range("Q12").Value = 6.916133360422501
$ws.Range("R12").Value = 41.496800162535
$ws.Range("S12").Value = 0.01470969349706466
$ws.Range("T12").Value = 0.0111772291826922

$ws.Range("G13").Value = 1.1448225
$ws.Range("H13").Value = 2.289645
$ws.Range("I13").Value = 0.08086486456196039
$ws.Range("J13").Value = 0.05540330324572383
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 9.354491666666666
$ws.Range("N13").Value = 28.063475
$ws.Range("O13").Value = 0.2816688040634135
$ws.Range("P13").Value = 0.3123873844340412
$ws.Range("Q13").Value = 10.7092325360625
$ws.Range("R13").Value = 64.25539521637501
$ws.Range("S13").Value = 0.02277710969191729
$ws.Range("T13").Value = 0.0173072929899377

$ws.Range("G14").Value = 1.1448225
$ws.Range("H14").Value = 2.289645
$ws.Range("I14").Value = 0.08086486456196039
$ws.Range("J14").Value = 0.05540330324572383
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 5.096198666666667
$ws.Range("N14").Value = 15.288596
$ws.Range("O14").Value = 0.1534492984610312
$ws.Range("P14").Value = 0.1701843594247949
$ws.Range("Q14").Value = 5.834242898070001
$ws.Range("R14").Value = 35.00545738842001
$ws.Range("S14").Value = 0.01240865673717912
$ws.Range("T14").Value = 0.00942877567289117

$ws.Range("G15").Value = 1.1448225
$ws.Range("H15").Value = 2.289645
$ws.Range("I15").Value = 0.08086486456196039
$ws.Range("J15").Value = 0.05540330324572383
$ws.Range("K15").Value = 2
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 9.797388999999999
$ws.Range("N15").Value = 19.594778
$ws.Range("O15").Value = 0.2950046823396649
$ws.Range("P15").Value = 0.2181184421382489
$ws.Range("Q15").Value = 11.2162713684525
$ws.Range("R15").Value = 44.86508547381
$ws.Range("S15").Value = 0.02385551368254115
$ws.Range("T15").Value = 0.01208448219327027

$ws.Range("G16").Value = 1.1448225
$ws.Range("H16").Value = 2.289645
$ws.Range("I16").Value = 0.08086486456196039
$ws.Range("J16").Value = 0.05540330324572383
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 2.921654
$ws.Range("N16").Value = 8.764962000000001
$ws.Range("O16").Value = 0.08797258230498059
$ws.Range("P16").Value = 0.09756680360660122
$ws.Range("Q16").Value = 3.344775236415
$ws.Range("R16").Value = 20.06865141849
$ws.Range("S16").Value = 0.03930929545995401
$ws.Range("T16").Value = 0.04480402713221842

